$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the 5 oldest data rows (rows 2-6), shifting everything else up
$ws.Range("2:6").Delete()

# Update selection to B2
$ws.Range("B2").Select()

# Update the SpreadsheetBuilder_2 defined name range
$wb.Names.Item("SpreadsheetBuilder_2").RefersTo = "=Sheet1!`$F`$2:`$G`$2"
